$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 4")
$ws.Activate()

# --- Row 6: new "Worked out" note in column I -------------------------------
# (written later, after the new "Lesson 4.8"/"Lesson 4.9" strings exist, so
#  the shared-string table picks up the same ordering as the target file)

# --- Row 11: add "P" marker in E11 ------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E11").Value = "P"

# --- Row 12: change Office Hours/6:30-7:00pm -> Class/7:00-9:00pm ----------
$ws.Range("E2").Copy()
$ws.Range("E12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E12").Value = "P"
$ws.Range("D12").Value = "Class"
$ws.Range("F12").Value = "7:00 pm - 9:00 pm"

# --- Row 13: this row's old content (Class / 7:00-9:00pm) is superseded by
#     the new row 17 below, so wipe it out entirely ---------------------------
$ws.Range("D13:F13").Clear()

# --- Row 14 (new): Wednesday ------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Wednesday"

$ws.Range("C2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "3/30/2022"

$ws2 = $wb.Worksheets.Item("Week 3")
$ws2.Range("G13").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "Lesson 4.8"

$ws.Range("E2").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Value = "P"

$ws.Range("F14").Value = "7:00 pm - 9:00 pm"
$ws.Range("G14").Value = 2

# --- Row 16 (new): Thursday -------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Thursday"

$ws.Range("C2").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "3/31/2022"

$ws.Range("D16").Value = "Lesson 4.9"

$ws.Range("E2").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "P"

$ws.Range("F16").Value = "5:30 am - 6:30 am"
$ws.Range("G16").Value = 1

# --- Row 17 (new): continuation line (date/P columns blank but formatted) --
$ws.Range("C2").Copy()
$ws.Range("C17").PasteSpecial(-4122)

$ws.Range("D17").Value = "Class"

$ws.Range("E2").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("F17").Value = "7:00 pm - 9:00 pm"

# --- Row 19 (new): Friday ----------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "Friday"

$ws.Range("C2").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "4/1/2022"

$ws.Range("F19").Value = "5:30 am - 6:30 am"
$ws.Range("G19").Value = 2

# --- Row 20 (new) -------------------------------------------------------------
$ws.Range("F20").Value = "5:30 pm - 6:30 pm"

# --- "Worked out" notes in column I (added after Lesson 4.8 / Lesson 4.9 so
#     shared-string indices land in the same order as the target workbook) --
$ws.Range("I6").Value = "Worked out"
$ws.Range("I11").Value = "Worked out"
$ws.Range("I14").Value = "Worked out"

# --- Update the Total Hours sum range + recalc ------------------------------
$ws.Range("G32").Formula = "=SUM(G2:G31)"

# --- Selection moves to D19, matching the author's final cursor position ---
$ws.Range("D19").Select()
